$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A (group id) and column B (count) values as per the diff.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 21355

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 13469

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 7724

$ws.Range("B5").Value = 4491
